$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Route")

# Update report date from 1/10/24 to 2/10/24
$ws.Range("L3").Value = (Get-Date -Year 2024 -Month 2 -Day 10 -Hour 0 -Minute 0 -Second 0).Date

# Update Unit Cost values for RSO 02, RSO 04, RSO 05 (route cost rows)
$ws.Range("D7").Value = 200
$ws.Range("D9").Value = 250
$ws.Range("D10").Value = 150

$wb.Save()
